$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "39.131.33"
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("D3").Value = "2.199.72"
$ws.Range("E3").Value = "  -6.98%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "295.60"
$ws.Range("E5").Value = "  -4.77%  "
$ws.Range("D6").Value = "82.30"
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("E7").Value = "  -3.60%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -5.81%  "
$ws.Range("E10").Value = "  -7.87%  "
$ws.Range("E11").Value = "  -5.10%  "
$ws.Range("E12").Value = "  -10.49%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "2.543.23"
$ws.Range("E14").Value = "  -6.87%  "
$ws.Range("D15").Value = "6.22"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").Value = "14.04"
$ws.Range("E16").Value = "  -6.53%  "
$ws.Range("D17").Value = "2.206.90"
$ws.Range("E17").Value = "  -5.57%  "
$ws.Range("D18").Value = "0.711"
$ws.Range("E18").Value = "  -5.82%  "
$ws.Range("D19").Value = "39.054.87"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").Value = "0.0₃0868"
$ws.Range("E20").Value = "  -4.35%  "
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -7.24%  "
$ws.Range("D22").Value = "64.47"
$ws.Range("E22").Value = "  -5.66%  "
$ws.Range("D23").Value = "10.20"
$ws.Range("E23").Value = "  -4.59%  "
$ws.Range("D24").Value = "227.33"
$ws.Range("E24").Value = "  -3.32%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D27").Value = "1.77"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "22.41"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "9.05"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "149.23"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").Value = "31.62"
$ws.Range("E32").Value = "  -6.98%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -7.33%  "
$ws.Range("D35").Value = "0.0690"
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("D38").Value = "0.0955"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("D39").Value = "2.62"
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("D40").Value = "14.88"
$ws.Range("E40").Value = "  -7.26%  "
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -4.11%  "
$ws.Range("D43").Value = "1.906.37"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").Value = "0.0258"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("E45").Value = "  -15.37%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.88"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.61"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").Value = "15.76"
$ws.Range("E48").Value = "  -11.03%  "
$ws.Range("D49").Value = "2.410.44"
$ws.Range("E49").Value = "  -7.26%  "
$ws.Range("D50").Value = "69.97"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").Value = "85.91"
$ws.Range("E51").Value = "  -7.31%  "
